# egg_notes_EGK.xlsx update: sort existing log by image number, append new
# entries (R script run added more processed images), and tidy up column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sort the existing data block (A2:B49) ascending by column A ---------
# (mirrors Excel's Data > Sort applied to the image-number / issue table)
$ws.Range("A2:B49").Sort($ws.Range("A2:A49"))

# --- 2. Append the new image filenames in the order they were typed --------
# (this also controls the order new entries land in the shared-string table)
$ws.Range("A51").Value = "IMG_1880.JPG"
$ws.Range("A52").Value = "IMG_2049.JPG"
$ws.Range("A53").Value = "IMG_2124.JPG"
$ws.Range("A54").Value = "IMG_2137.JPG"
$ws.Range("A55").Value = "IMG_2139.JPG"
$ws.Range("A56").Value = "IMG_2142.JPG"
$ws.Range("A57").Value = "IMG_2178.JPG"
$ws.Range("A58").Value = "IMG_2525.JPG"
$ws.Range("A59").Value = "IMG_2528.JPG"
$ws.Range("A60").Value = "IMG_2822.JPG"
$ws.Range("A61").Value = "IMG_2921.JPG"
$ws.Range("A62").Value = "IMG_3035.JPG"
$ws.Range("A63").Value = "IMG_3056.JPG"
$ws.Range("A64").Value = "IMG_3085.JPG"
$ws.Range("A65").Value = "IMG_3442.JPG"

# --- 3. A new header/flag row was inserted above the new list --------------
$ws.Range("A50").Value = "Bad_thresh"

# --- 4. Resize column A to fit the new (longer) filenames ------------------
$ws.Columns("A").ColumnWidth = 12.1640625

# --- 5. Restore the on-screen selection / scroll position left by the edit -
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C53").Select() | Out-Null
